# Fixing Legs For 1,2, and 3
#
# The "L3" row of raw calibration readings (row 9) is corrected:
#   D9 (alpha 0 combo B reading): 1350 -> 1300
#   G9 (alpha 90 combo A reading): 900 -> 850
# H9's shared formula (=D9-2*(D9-G9)) recalculates automatically from
# 450 to 400 once the inputs above are written.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D9").Value = 1300
$ws.Range("G9").Value = 850

# Reflect the author's final scroll position/selection on the sheet.
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1
$null = $ws.Range("G11").Select()
